# Insert a new data row at row 83 (pushing the existing rows 83-178 down to
# 84-179) and populate it with a new Albahaca price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 83..178 down by one (xlShiftDown = -4121), creating a blank row 83.
$ws.Rows("83:83").Insert(-4121)

# Populate the newly inserted row 83 with the new record.
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44159
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 100112052
$ws.Range("G83").Value = "Albahaca"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 80
$ws.Range("K83").Value = 6000
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = 6000
$ws.Range("N83").Value = "`$/docena de matas"
$ws.Range("O83").Value = "Provincia de Quillota"
$ws.Range("P83").Value = 1000
$ws.Range("Q83").Value = 6
$ws.Range("R83").Value = "Hortaliza"
